# Add an "InvalidLogin" worksheet (after "ValidLogin") with invalid
# credentials, to verify the application's error-message scenario.

$wb = $excel.ActiveWorkbook

# Create the new sheet and give it its final name first.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "InvalidLogin"

# Re-fetch "ValidLogin" fresh (by name) right before using it so we don't
# rely on a reference captured before the sheet collection changed order.
$newSheet.Move($null, $wb.Worksheets.Item("ValidLogin"))

# Populate the new sheet with the invalid-login test data.
$ws2 = $wb.Worksheets.Item("InvalidLogin")
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "abc"
$ws2.Range("B2").Value = "xyz"

# Select the whole first row on the original sheet (as left by the author).
$ws1 = $wb.Worksheets.Item("ValidLogin")
$null = $ws1.Rows("1:1").Select()

# Make the new sheet the active / selected tab, with B2 selected.
$ws2 = $wb.Worksheets.Item("InvalidLogin")
$null = $ws2.Activate()
$null = $ws2.Range("B2").Select()
